# "remove old average code"
# Clears the stale average/percentage figures that used to live in
# columns E:J for rows 13-54, leaving the cells blank but keeping the
# same "empty data" formatting already used by the neighbouring K:P
# columns on each of those rows (cellXf style index 30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 13
$lastRow = 54

# K13 (and every K cell down through K54) already carries the blank-cell
# style used across the sheet for "no data" columns. Copying its format
# onto E:J reproduces that exact style (rather than fabricating a new
# style entry), then ClearContents wipes the leftover numeric values
# without touching the freshly applied formatting.
$styleSource = $ws.Range("K" + $firstRow)
$target = $ws.Range("E" + $firstRow + ":J" + $lastRow)

$styleSource.Copy()
$target.PasteSpecial(-4122)
$target.ClearContents()

$excel.CutCopyMode = $false
